# Worker List.xlsx edit:
# Remove the stray "jess" row (row 35: jess / 1219 / 70626308 / N521D5060018)
# from the WorkerList table on Sheet1. This shifts everything below it up
# by one row, shrinks the table/autofilter/used range by one row, and
# (since those two shared strings were only referenced by that row) the
# shared-strings table also loses those two now-unused entries when Excel
# rewrites the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 35 (shifts rows 36.. up by one).
$ws.Rows(35).Delete()

# Fix up the hidden _xlnm._FilterDatabase defined name so it tracks the
# table's new (one-row-shorter) extent.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$4:`$A`$51"
    }
}

# Re-apply the existing (already-alphabetical, so data doesn't actually
# move) sort so the sheet's remembered sort-range shrinks along with the
# deleted row, instead of staying stale at its old extent.
$sortRange = $ws.Range("A4:C47")
$sortKey = $ws.Range("A1")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
[void]$ws.Sort.Apply()

# Restore the last-used selection/cursor position recorded in the sheet.
[void]$ws.Range("A39").Select()
